# Wheel Bill of Materials Cost Weight - edit script
# Implements: new "per car" columns (Packs to Order per Car / Car Price),
# new informational note (B1), new Metalbits/Maplin subtotal rows (7 & 16),
# per-unit pricing columns (Pack Price/Pack Qty -> Unit Price, per Unit of
# Measure) driving the existing Wheel-Qty/Wheel-Price math, and assorted
# formatting (new £-with-more-decimals number formats, size-14 subtotal
# font).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Number format codes reused from the existing workbook (Currency cell
# style variants) plus the two new higher-precision £ formats.
# ---------------------------------------------------------------------------
$fmtCur2   = '_-"£"* #,##0.00_-;\-"£"* #,##0.00_-;_-"£"* "-"??_-;_-@_-'
$fmtCur4   = '_-"£"* #,##0.0000_-;\-"£"* #,##0.0000_-;_-"£"* "-"??_-;_-@_-'
$fmtCur5   = '_-"£"* #,##0.00000_-;\-"£"* #,##0.00000_-;_-"£"* "-"??_-;_-@_-'
$fmtComma0 = '_-* #,##0_-;\-* #,##0_-;_-* "-"??_-;_-@_-'

# ---------------------------------------------------------------------------
# Step 1: make room for the new columns.
#   old E..H (Unit Price/Wheel Qty/Unit of Measure/Packs to Order) shift
#   right by one (a new "Unit of Measure" column is inserted right before
#   them, at old column E) and old J..N (Car Price/Supplier/Shipping/Date of
#   Order/Ship From) shift right by three more (two new columns - "Packs to
#   Order per Car" and a blank spacer - are inserted before old column J).
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).Insert()
$ws.Columns.Item(11).Insert()
$ws.Columns.Item(11).Insert()

# ---------------------------------------------------------------------------
# Step 2: header row (row 1). New B1 note, new E1/H1 "Unit of Measure"
# columns, new K1 "Packs to Order per Car" column, everything else keeps
# its previous text (now one or three columns further right).
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "This spreadsheet excludes shipping, as that decreases with volume a lot.  It includes VAT. Or should it?"
$ws.Range("E1").Value = "Unit of Measure"
$ws.Range("H1").Value = "Unit of Measure"
$ws.Range("K1").Value = "Packs to Order per Car"

$ws.Range("B1:Q1").WrapText = $true
$ws.Range("K1").Font.Bold = $false
$ws.Range("L1").Clear()
$ws.Rows.Item(1).RowHeight = 45

Write-Output "row1 done"

# ---------------------------------------------------------------------------
# Step 3: Metalbits rows (4-6) - add the new Pack Price / Pack Qty -> Unit
# Price chain (C/D/F), the doubled Wheel-Qty formula (G), the new
# Packs-to-Order / per-car columns (I/J/K/M) and the "mm" unit labels
# (E/H/L). Row 4 = 12mm pads, row 5 = 10mm side pads, row 6 = 6mm rod.
# ---------------------------------------------------------------------------

# Row 4 - 12mm x 12mm x 12mm Pads
$ws.Range("C4").Formula = "=3.32/2*1.2"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "mm"
$ws.Range("F4").Formula = "=C4/D4"
$ws.Range("G4").Formula = "=48*13*2"
$ws.Range("H4").Value = "mm"
$ws.Range("I4").Value = 2
$ws.Range("J4").Formula = "=I4*F4"
$ws.Range("K4").Formula = "=4*G4"
$ws.Range("L4").Value = "mm"
$ws.Range("M4").Formula = "=5*F4"

$ws.Range("C4").NumberFormat = $fmtCur2
$ws.Range("F4").NumberFormat = $fmtCur2
$ws.Range("J4").NumberFormat = $fmtCur2
$ws.Range("J4").Font.Bold = $true
$ws.Range("M4").NumberFormat = $fmtCur2
$ws.Range("M4").Font.Bold = $true
$ws.Range("K4:L4").Font.Bold = $false
$ws.Range("K4:L4").NumberFormat = "General"

# Row 5 - 10mm x 3mm Side Pads & Top Circuit
$ws.Range("B5").Value = "10mm x 3mm Side Pads & Top Circuit, cuts of 3mm.  Order qty is 1m"
$ws.Range("C5").Formula = "=0.59*1.2"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "mm"
$ws.Range("F5").Formula = "=C5/D5"
$ws.Range("G5").Formula = "=48*(15+15+121)"
$ws.Range("H5").Value = "mm"
$ws.Range("I5").Value = 8
$ws.Range("J5").Formula = "=I5*F5"
$ws.Range("K5").Formula = "=4*G5"
$ws.Range("L5").Value = "mm"
$ws.Range("M5").Formula = "=29*F5"

$ws.Range("B5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 30
$ws.Range("C5").NumberFormat = $fmtCur4
$ws.Range("F5").NumberFormat = $fmtCur2
$ws.Range("J5").NumberFormat = $fmtCur2
$ws.Range("J5").Font.Bold = $true
$ws.Range("M5").NumberFormat = $fmtCur2
$ws.Range("M5").Font.Bold = $true
$ws.Range("K5:L5").Font.Bold = $false
$ws.Range("K5:L5").NumberFormat = "General"

# Row 6 - 6mm rod x 100mm
$ws.Range("C6").Formula = "=0.5*1.2"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = "mm"
$ws.Range("F6").Formula = "=C6/D6"
$ws.Range("H6").Value = "mm"
$ws.Range("I6").Value = 5
$ws.Range("J6").Formula = "=I6*F6"
$ws.Range("K6").Formula = "=4*G6"
$ws.Range("L6").Value = "mm"
$ws.Range("M6").Formula = "=20*F6"

$ws.Range("C6").NumberFormat = $fmtCur2
$ws.Range("F6").NumberFormat = $fmtCur2
$ws.Range("J6").NumberFormat = $fmtCur2
$ws.Range("J6").Font.Bold = $true
$ws.Range("M6").NumberFormat = $fmtCur2
$ws.Range("M6").Font.Bold = $true
$ws.Range("K6:L6").Font.Bold = $false
$ws.Range("K6:L6").NumberFormat = "General"

Write-Output "metalbits rows done"

# ---------------------------------------------------------------------------
# Step 4: Metalbits subtotal row (new row 7), size-14 font, bold label.
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = "Metalbits Subtotal for one wheel, one car"
$ws.Range("J7").Formula = "=SUM(J4:J6)"
$ws.Range("M7").Formula = "=SUM(M4:M6)"

$ws.Range("B7:O7").Font.Size = 14
$ws.Range("B7").Font.Bold = $true
$ws.Range("J7").Font.Bold = $true
$ws.Range("M7").Font.Bold = $true
$ws.Range("J7").NumberFormat = $fmtCur2
$ws.Range("M7").NumberFormat = $fmtCur2
$ws.Range("C7").NumberFormat = $fmtCur2
$ws.Range("O7").NumberFormat = $fmtCur2
$ws.Rows.Item(7).RowHeight = 18

Write-Output "metalbits subtotal row done"
